# Apply numeric corrections to leve profit calculation sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 53585940
$ws.Range("I74").Value = 107153496
$ws.Range("K74").Value = 107153496
$ws.Range("M74").Value = -107152560

$ws.Range("H77").Value = 53585940
$ws.Range("I77").Value = 107153496
$ws.Range("K77").Value = 535767480
$ws.Range("M77").Value = -535762800

$ws.Range("H111").Value = 5683310
$ws.Range("I111").Value = 10416928
$ws.Range("J111").Value = 2968.9
$ws.Range("K111").Value = 31250784
$ws.Range("L111").Value = 8906.700000000001
$ws.Range("M111").Value = -31247717
$ws.Range("N111").Value = -15040.7

$ws.Range("H121").Value = 909.1667
$ws.Range("J121").Value = 1001
$ws.Range("L121").Value = 3003
$ws.Range("N121").Value = -6497

$ws.Range("H138").Value = 2218.6
$ws.Range("J138").Value = 2607.5667
$ws.Range("L138").Value = 7822.7001
$ws.Range("N138").Value = -18102.7001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4210.875
$ws.Range("I45").Value = 2489.923
$ws.Range("J45").Value = 6244.727
$ws.Range("K45").Value = 2489.923
$ws.Range("L45").Value = 6244.727
$ws.Range("M45").Value = -2112.923
$ws.Range("N45").Value = -6998.727

$ws.Range("H61").Value = 26318408
$ws.Range("I61").Value = 1892.069
$ws.Range("J61").Value = 111116070
$ws.Range("K61").Value = 1892.069
$ws.Range("L61").Value = 111116070
$ws.Range("M61").Value = -1680.069
$ws.Range("N61").Value = -111116494

$ws.Range("H74").Value = 29066.447
$ws.Range("I74").Value = 43055.625
$ws.Range("K74").Value = 43055.625
$ws.Range("M74").Value = -42181.625

$ws.Range("H77").Value = 29066.447
$ws.Range("I77").Value = 43055.625
$ws.Range("K77").Value = 215278.125
$ws.Range("M77").Value = -210910.125

$ws.Range("H119").Value = 48734.25
$ws.Range("J119").Value = 48734.25
$ws.Range("L119").Value = 48734.25
$ws.Range("N119").Value = -58410.25

$ws.Range("H136").Value = 26318408
$ws.Range("I136").Value = 1892.069
$ws.Range("J136").Value = 111116070
$ws.Range("K136").Value = 5676.207
$ws.Range("L136").Value = 333348210
$ws.Range("M136").Value = -3126.207
$ws.Range("N136").Value = -333353310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 32818.453
$ws.Range("I86").Value = 47563.684
$ws.Range("J86").Value = 3328
$ws.Range("K86").Value = 47563.684
$ws.Range("L86").Value = 3328
$ws.Range("M86").Value = -46440.684
$ws.Range("N86").Value = -5574

$ws.Range("H89").Value = 32818.453
$ws.Range("I89").Value = 47563.684
$ws.Range("J89").Value = 3328
$ws.Range("K89").Value = 237818.42
$ws.Range("L89").Value = 16640
$ws.Range("M89").Value = -232202.42
$ws.Range("N89").Value = -27872

$ws.Range("H94").Value = 2513.1516
$ws.Range("I94").Value = 1183.5217
$ws.Range("K94").Value = 1183.5217
$ws.Range("M94").Value = -732.5217

$ws.Range("H107").Value = 13893022
$ws.Range("I107").Value = 17859600
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 17859600
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = -17857680
$ws.Range("N107").Value = -13840

$ws.Range("H134").Value = 9264468
$ws.Range("I134").Value = 35715696
$ws.Range("J134").Value = 6538.35
$ws.Range("K134").Value = 107147088
$ws.Range("L134").Value = 19615.05
$ws.Range("M134").Value = -107144553
$ws.Range("N134").Value = -24685.05

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5617.244
$ws.Range("I31").Value = 1267.909
$ws.Range("K31").Value = 1267.909
$ws.Range("M31").Value = -972.9090000000001

$ws.Range("H34").Value = 5617.244
$ws.Range("I34").Value = 1267.909
$ws.Range("K34").Value = 1267.909
$ws.Range("M34").Value = -1065.909

$ws.Range("H132").Value = 3879.7046
$ws.Range("I132").Value = 2345.3225
$ws.Range("J132").Value = 7538.615
$ws.Range("K132").Value = 7035.967500000001
$ws.Range("L132").Value = 22615.845
$ws.Range("M132").Value = -4505.967500000001
$ws.Range("N132").Value = -27675.845

$ws.Range("H134").Value = 2737.4
$ws.Range("I134").Value = 1377.4517
$ws.Range("K134").Value = 4132.355100000001
$ws.Range("M134").Value = -1597.355100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 777
$ws.Range("I86").Value = 375
$ws.Range("K86").Value = 1125
$ws.Range("M86").Value = 61

$ws.Range("H89").Value = 777
$ws.Range("I89").Value = 375
$ws.Range("K89").Value = 3375
$ws.Range("M89").Value = 2553

$ws.Range("H113").Value = 7464.7334
$ws.Range("I113").Value = 1897.1666
$ws.Range("J113").Value = 11176.444
$ws.Range("K113").Value = 5691.4998
$ws.Range("L113").Value = 33529.33199999999
$ws.Range("M113").Value = -3521.4998
$ws.Range("N113").Value = -37869.33199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2318.7144
$ws.Range("J80").Value = 2219.4
$ws.Range("L80").Value = 2219.4
$ws.Range("N80").Value = -4215.4

$ws.Range("H83").Value = 2318.7144
$ws.Range("J83").Value = 2219.4
$ws.Range("L83").Value = 11097
$ws.Range("N83").Value = -21081

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 886.9655
$ws.Range("I22").Value = 278.3125
$ws.Range("J22").Value = 1636.0769
$ws.Range("K22").Value = 278.3125
$ws.Range("L22").Value = 1636.0769
$ws.Range("M22").Value = 16.6875
$ws.Range("N22").Value = -2226.0769

$ws.Range("H27").Value = 886.9655
$ws.Range("I27").Value = 278.3125
$ws.Range("J27").Value = 1636.0769
$ws.Range("K27").Value = 278.3125
$ws.Range("L27").Value = 1636.0769
$ws.Range("M27").Value = -171.3125
$ws.Range("N27").Value = -1850.0769

$ws.Range("H38").Value = 21714.5
$ws.Range("J38").Value = 21714.5
$ws.Range("L38").Value = 21714.5
$ws.Range("N38").Value = -22534.5

$ws.Range("H119").Value = 56134
$ws.Range("J119").Value = 56134
$ws.Range("L119").Value = 56134
$ws.Range("N119").Value = -65810

$ws.Range("H136").Value = 8152.183
$ws.Range("J136").Value = 12985.934
$ws.Range("L136").Value = 38957.802
$ws.Range("N136").Value = -44057.802

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 56133.5
$ws.Range("J119").Value = 56133.5
$ws.Range("L119").Value = 56133.5
$ws.Range("N119").Value = -65809.5
